$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.060.82"
$ws.Range("E2").Value = "  +2.59%  "
$ws.Range("D3").Value = "2.471.91"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'567.66"
$ws.Range("E5").Value = "  +1.92%  "
$ws.Range("D6").Value = "'166.92"
$ws.Range("E6").Value = "  +2.84%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "'0.513"
$ws.Range("E8").Value = "  +0.74%  "
$ws.Range("D9").Value = "'0.177"
$ws.Range("E9").Value = "  +13.87%  "
$ws.Range("E10").Value = "  -1.25%  "
$ws.Range("D11").Value = "'0.335"
$ws.Range("E11").Value = "  +2.28%  "
$ws.Range("D12").Value = "'4.67"
$ws.Range("E12").Value = "  -2.93%  "
$ws.Range("D13").Value = "'0.0000183"
$ws.Range("E13").Value = "  +8.46%  "
$ws.Range("D14").Value = "69.882.76"
$ws.Range("E14").Value = "  +2.44%  "
$ws.Range("D15").Value = "2.925.00"
$ws.Range("E15").Value = "  +0.66%  "
$ws.Range("D16").Value = "'24.08"
$ws.Range("E16").Value = "  +3.30%  "
$ws.Range("D17").Value = "2.471.10"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").Value = "'10.84"
$ws.Range("E18").Value = "  +3.57%  "
$ws.Range("D19").Value = "'343.33"
$ws.Range("E19").Value = "  +2.02%  "
$ws.Range("D20").Value = "'7.18"
$ws.Range("E20").Value = "  +4.25%  "
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").Value = "'3.90"
$ws.Range("E21").Value = "  +3.18%  "
$ws.Range("B22").Value = "SuiNetwork"
$ws.Range("C22").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D22").Value = "'2.04"
$ws.Range("E22").Value = "  +8.95%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "'66.69"
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("D25").Value = "'3.91"
$ws.Range("E25").Value = "  +6.39%  "
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "2.594.98"
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").Value = "'8.55"
$ws.Range("E27").Value = "  +5.36%  "
$ws.Range("D28").Value = "'0.990"
$ws.Range("E28").Value = "  -0.87%  "
$ws.Range("D29").Value = "0.0₃0854"
$ws.Range("E29").Value = "  +4.93%  "
$ws.Range("D30").Value = "'7.34"
$ws.Range("E30").Value = "  +2.16%  "
$ws.Range("D31").Value = "'1.25"
$ws.Range("E31").Value = "  +10.05%  "
$ws.Range("D32").Value = "'447.86"
$ws.Range("E32").Value = "  +5.20%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("D34").Value = "'1.63"
$ws.Range("E34").Value = "  +1.28%  "
$ws.Range("D35").Value = "'160.98"
$ws.Range("E35").Value = "  +0.36%  "
$ws.Range("D36").Value = "'19.07"
$ws.Range("E36").Value = "  +0.48%  "
$ws.Range("B37").Value = "USDe"
$ws.Range("C37").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.110"
$ws.Range("E38").Value = "  +4.15%  "
$ws.Range("E39").Value = "  +2.56%  "
$ws.Range("D40").Value = "'0.307"
$ws.Range("E40").Value = "  +3.96%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.56"
$ws.Range("E41").Value = "  +6.09%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").Value = "'4.49"
$ws.Range("E42").Value = "  +2.86%  "
$ws.Range("E43").Value = "  +3.22%  "
$ws.Range("D44").Value = "'2.14"
$ws.Range("E44").Value = "  +6.69%  "
$ws.Range("D45").Value = "'3.41"
$ws.Range("E45").Value = "  +1.65%  "
$ws.Range("D46").Value = "'132.98"
$ws.Range("E46").Value = "  +2.55%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.0726"
$ws.Range("E47").Value = "  +1.18%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'0.492"
$ws.Range("E48").Value = "  +2.44%  "
$ws.Range("D49").Value = "'0.565"
$ws.Range("E49").Value = "  +0.73%  "
$ws.Range("D50").Value = "'0.0931"
$ws.Range("E50").Value = "  +1.54%  "
$ws.Range("E51").Value = "  +3.01%  "

# Strip the quote-prefix style marker left behind on text-forced numeric-looking cells
$ws.Range("D4").ClearFormats()
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D12").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D18").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D50").ClearFormats()
